$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:F10 (ME, MAE, MSE, RMSE, SE) per row (Q0..Q8)
$data = @{
    2 = @(0.2033919902778603, 0.4364632531330045, 0.3742668131629021, 0.6117734982515196, 0.5987537801616326)
    3 = @(0.2628942463139419, 0.4333144524565598, 0.2926323150708423, 0.5409550028152456, 0.4920828262122367)
    4 = @(0.2853914118626144, 0.4746291375408727, 0.3229524026090269, 0.5682890132749593, 0.5132826384015559)
    5 = @(0.2197630776494257, 0.4033531803712247, 0.234430775428361,  0.4841805194639299, 0.4524913939993277)
    6 = @(0.1926414270218137, 0.3605436531697729, 0.1752150039837381, 0.4185869132972722, 0.3917259310889259)
    7 = @(0.2387356972361521, 0.2818843305516119, 0.1725877670096401, 0.4154368387729236, 0.3606135925185842)
    8 = @(0.3398310944915218, 0.3398310944915218, 0.2116720379002942, 0.4600782954023089, 0.3397414283545448)
    9 = @(0.4537447083867415, 0.4537447083867415, 0.2538588652253824, 0.5038440882112069, 0.2682571662688996)
    10 = @(0.08161590864515353, 0.08161590864515353, 0.006661156543974047, 0.08161590864515353)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $col = 2  # Column B
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col++
    }
}
